$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 statistics update (May 2025 review 2)
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 15.85638515958219
$ws.Range("F4").Value = 1.035985783041232
$ws.Range("G4").Value = 14.89943443887979
$ws.Range("H4").Value = 15.40321241818091
$ws.Range("I4").Value = 15.59840741315882
$ws.Range("J4").Value = 16.05158015456011
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 19.43074970151516
$ws.Range("N4").Value = 2.327761773041207
$ws.Range("P4").Value = 18.14331576257711
$ws.Range("R4").Value = 19.67246997689321
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 22.78301885253688
$ws.Range("V4").Value = 3.060001819657683
$ws.Range("W4").Value = 20.14183710066337
$ws.Range("X4").Value = 20.99018518967267
$ws.Range("Y4").Value = 21.93562612555283
$ws.Range("Z4").Value = 23.72845978841703
$ws.Range("AR4").Value = 4
$ws.Range("AS4").Value = 6.92663369295469
$ws.Range("AT4").Value = 2.055932201174247
$ws.Range("AU4").Value = 5.242402661783577
$ws.Range("AV4").Value = 5.546069268675453
$ws.Range("AW4").Value = 6.337218712394003
$ws.Range("AX4").Value = 7.717783136673241
$ws.Range("AZ4").Value = 4
$ws.Range("BA4").Value = 3.574364541932967
$ws.Range("BB4").Value = 1.381235801940727
$ws.Range("BD4").Value = 2.873553838038495
$ws.Range("BE4").Value = 3.149749612874622
$ws.Range("BF4").Value = 3.850560316769094
$ws.Range("BH4").Value = 4
$ws.Range("BI4").Value = 2.899977208941098
$ws.Range("BJ4").Value = 2.055932201174247
$ws.Range("BK4").Value = 1.215746177769985
$ws.Range("BL4").Value = 1.51941278466186
$ws.Range("BM4").Value = 2.310562228380411
$ws.Range("BN4").Value = 3.691126652659648
$ws.Range("BP4").Value = 4
$ws.Range("BQ4").Value = 3.312889829991862
$ws.Range("BR4").Value = 1.381235801940727
$ws.Range("BT4").Value = 2.61207912609739
$ws.Range("BU4").Value = 2.888274900933517
$ws.Range("BV4").Value = 3.589085604827989
$ws.Range("BX4").Value = 4
$ws.Range("BY4").Value = 0.2247786069507456
$ws.Range("BZ4").Value = 0.1870609715601037
$ws.Range("CB4").Value = 0.09832174606253023
$ws.Range("CC4").Value = 0.225074950136276
$ws.Range("CD4").Value = 0.3515318110244913
$ws.Range("CE4").Value = 0.4305503368806072
$ws.Range("CF4").Value = 4
$ws.Range("CG4").Value = 0.1299733728859163
$ws.Range("CH4").Value = 0.08184749272893087
$ws.Range("CJ4").Value = 0.09912146140571043
$ws.Range("CK4").Value = 0.1355791086939757
$ws.Range("CL4").Value = 0.1664310201741815
$ws.Range("CN4").Value = 0.09353048578005187
$ws.Range("CO4").Value = 0.04092374636446543
